$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for the three new vehicles (AEC Regent I, AEC Regent II, AEC Routemaster)
$vehicles = @(
    @{ Row = 24; Name = "AEC Regent I";   Year = 1929; F = 29; G = 50 },
    @{ Row = 25; Name = "AEC Regent II";  Year = 1945; F = 38; G = 57 },
    @{ Row = 26; Name = "AEC Routemaster"; Year = 1954; F = 40; G = 72 }
)

foreach ($v in $vehicles) {
    $r = $v.Row

    $ws.Cells.Item($r, 1).Value = $v.Name
    $ws.Cells.Item($r, 2).Value = $v.Year
    $ws.Cells.Item($r, 3).Value = 1
    $ws.Cells.Item($r, 4).Value = "Bus"

    $ws.Cells.Item($r, 5).Formula = "=IF(B$r > 1900, ((B$r-1900)*10)+400+C$r, ((B$r-1730)*2)+C$r)+VLOOKUP(D$r,'ID Scheme'!`$A`$2:`$B`$4,2)"

    $ws.Cells.Item($r, 6).Value = $v.F
    $ws.Cells.Item($r, 7).Value = $v.G

    $ws.Cells.Item($r, 8).Formula = "=SQRT(F$r*G$r)/`$B`$1"
    $ws.Cells.Item($r, 8).NumberFormat = $ws.Cells.Item(23, 8).NumberFormat

    $ws.Cells.Item($r, 9).Formula = "=H$r*0.9"
    $ws.Cells.Item($r, 9).NumberFormat = $ws.Cells.Item(23, 9).NumberFormat

    $ws.Cells.Item($r, 10).Value = "x"
    $ws.Cells.Item($r, 10).NumberFormat = $ws.Cells.Item(23, 10).NumberFormat
}

$ws.Rows.Item(8).RowHeight = 14.25

$ws.Range("A4").Select()
